$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.686274509803922
$ws.Range("C2").Value = 0.689716312056738
$ws.Range("D2").Value = 0.804147465437788
$ws.Range("E2").Value = 0.667359667359667
$ws.Range("F2").Value = 0.550151975683891

$ws.Range("B3").Value = 0.803455723542117
$ws.Range("C3").Value = 0.823333333333333
$ws.Range("D3").Value = 0.850678733031674
$ws.Range("E3").Value = 0.798850574712644
$ws.Range("F3").Value = 0.677653902084343

$ws.Range("B4").Value = 0.741721854304636
$ws.Range("C4").Value = 0.79020979020979
$ws.Range("D4").Value = 0.775229357798165
$ws.Range("E4").Value = 0.709090909090909
$ws.Range("F4").Value = 0.592648539778449
